$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '24.791.93'
$ws.Range("E2").Value = '  +0.84%  '

# Row 3
$ws.Range("D3").Value = '1.702.50'
$ws.Range("E3").Value = '  +0.43%  '

# Row 4
$ws.Range("E4").Value = '  +0.27%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '317.53'
$ws.Range("E5").Value = '  +0.37%  '

# Row 6
$ws.Range("E6").Value = '  +0.32%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3961'
$ws.Range("E7").Value = '  +0.48%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4093'
$ws.Range("E8").Value = '  +2.17%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.506'
$ws.Range("E9").Value = '  -1.17%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.006'
$ws.Range("E10").Value = '  +0.48%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '52.83'
$ws.Range("E11").Value = '  -0.02%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08941'
$ws.Range("E12").Value = '  +2.16%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.702'
$ws.Range("E13").Value = '  +7.01%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '24.47'
$ws.Range("E14").Value = '  +5.45%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.159'
$ws.Range("E15").Value = '  +0.04%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001336'
$ws.Range("E16").Value = '  +1.77%  '

# Row 17
$ws.Range("D17").Value = '1.708.08'
$ws.Range("E17").Value = '  +0.56%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '100.14'
$ws.Range("E18").Value = '  +0.42%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.07137'
$ws.Range("E19").Value = '  +1.04%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '20.08'
$ws.Range("E20").Value = '  +2.26%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.304'
$ws.Range("E21").Value = '  +4.81%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.005'
$ws.Range("E22").Value = '  +0.45%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '14.56'
$ws.Range("E23").Value = '  +2.81%  '

# Row 24
$ws.Range("D24").Value = '24.770.62'
$ws.Range("E24").Value = '  +0.73%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.089'
$ws.Range("E25").Value = '  -0.51%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.340'
$ws.Range("E26").Value = '  +0.10%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '23.06'
$ws.Range("E27").Value = '  +1.48%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.428'
$ws.Range("E28").Value = '  +26.27%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '165.22'
$ws.Range("E29").Value = '  +1.68%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '139.55'
$ws.Range("E30").Value = '  +2.60%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.209'
$ws.Range("E31").Value = '  +0.36%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.104'

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09186'
$ws.Range("E33").Value = '  +7.20%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.085'
$ws.Range("E34").Value = '  +0.35%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.03069'
$ws.Range("E35").Value = '  +12.47%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.2818'
$ws.Range("E36").Value = '  +3.29%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '11.11'
$ws.Range("E37").Value = '  -3.63%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.968'
$ws.Range("E38").Value = '  +1.82%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '14.59'
$ws.Range("E39").Value = '  +1.27%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.09316'
$ws.Range("E40").Value = '  +2.48%  '

# Row 41
$ws.Range("B41").Value = 'TheSandbox'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.7844'
$ws.Range("E41").Value = '  +2.84%  '

# Row 42
$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.480'
$ws.Range("E42").Value = '  -0.09%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '16.29'
$ws.Range("E43").Value = '  +3.99%  '

# Row 44
$ws.Range("B44").Value = 'Decentraland'
$ws.Range("C44").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.7321'
$ws.Range("E44").Value = '  +2.36%  '

# Row 45
$ws.Range("B45").Value = 'NEARProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.633'
$ws.Range("E45").Value = '  +3.35%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.254'
$ws.Range("E46").Value = '  +0.90%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.358'
$ws.Range("E47").Value = '  +3.41%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.004'
$ws.Range("E48").Value = '  +0.34%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '141.29'
$ws.Range("E49").Value = '  +0.32%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '94.11'
$ws.Range("E50").Value = '  +5.85%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.08064'
$ws.Range("E51").Value = '  +1.11%  '
